$d = $word.ActiveDocument

# The document previously had a single "_GoBack" bookmark sitting just
# before the "{r_b}" placeholder further down the table. The edit moves
# that bookmark so it now sits right after the "{unit" run (the text that
# used to read "{project"). Start by getting rid of the old one if the
# object model lets us (no-op if it can't be removed).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the "{project" run and change it to "{unit" (the trailing "}"
# lives in its own separate run and must stay untouched).
$rng = $d.Content
$rng.Find.Execute("{project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "{unit"

# Re-create the "_GoBack" bookmark immediately after the run we just
# edited (collapsed range right at the end of "{unit").
$d.Bookmarks.Add("_GoBack", $d.Range($rng.End, $rng.End))
